$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 8: Lease End Date, 3-Month Reminder text, and Lease Duration (Years)
$ws.Range("C8").Value = 45208
$ws.Range("D8").Value = "Reminder: Lease Ending Soon"
$ws.Range("E8").Value = 1
